$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replicate the date cell's style (numFmtId 14) from A4, then set the new value
$ws.Range("A4").Copy($ws.Range("A5"))
$ws.Range("A5").Value = 43809

$ws.Range("B5").Value = "Design for the login,register and menu on this project."
$ws.Range("C5").Value = "Design for the project"

$ws.Range("C5").Select()
